# disconnected charcoal production from heat fuel.
# Removes the "connections" row that links heat (simple_heat) as an
# inflow of biofuel/charcoal, outflow simple_charcoal/charcoal.
# This was row 24 on the "connections" worksheet; deleting it shifts the
# rows below (old 25, 26) up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("connections")

# Make sure we're deleting the right row (sanity: B24="heat", C24="simple_heat").
$ws.Rows.Item(24).EntireRow.Delete() | Out-Null

# Match the saved selection/active cell shown in the post-edit file.
$ws.Activate()
$ws.Range("B30").Select() | Out-Null
